# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Thu Jun 15 18:48:15 UTC 2023 with GitHub Actions"
# Column D (Price) and Column E (Volume(1h)) values are refreshed for every
# coin row; row 47/48 additionally swap their Coin/Link/Price/Volume content
# (Cronos <-> Aptos reorder upstream).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.137.96"
$ws.Range("E2").Value = "  -3.06%  "
$ws.Range("D3").Value = "1.651.46"
$ws.Range("E3").Value = "  -4.97%  "
$ws.Range("D4").Value = "'0.9981"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'234.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.10%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4776"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.09%  "
$ws.Range("D8").Value = "'0.2567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.83%  "
$ws.Range("D9").Value = "'0.06107"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "'0.07042"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("D11").Value = "1.642.67"
$ws.Range("E11").Value = "  -5.44%  "
$ws.Range("D12").Value = "'14.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.16%  "
$ws.Range("D13").Value = "'0.5762"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -12.33%  "
$ws.Range("D14").Value = "'4.314"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.76%  "
$ws.Range("D15").Value = "'73.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.56%  "
$ws.Range("D16").Value = "'0.9990"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'0.9991"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "25.129.59"
$ws.Range("E18").Value = "  -3.07%  "
$ws.Range("D19").Value = "'0.000006646"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").Value = "'11.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.70%  "
$ws.Range("D21").Value = "1.852.91"
$ws.Range("E21").Value = "  -5.37%  "
$ws.Range("D22").Value = "'4.323"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.02%  "
$ws.Range("D23").Value = "'8.475"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("D24").Value = "'5.242"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.34%  "
$ws.Range("D25").Value = "'134.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'14.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").Value = "'1.374"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.47%  "
$ws.Range("D28").Value = "'103.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").Value = "'1.645"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.14%  "
$ws.Range("D30").Value = "'3.924"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").Value = "'0.07608"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.33%  "
$ws.Range("D32").Value = "'3.541"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("D33").Value = "'0.9987"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'0.04292"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.25%  "
$ws.Range("D35").Value = "'2.576"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("D36").Value = "'0.9363"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.61%  "
$ws.Range("D37").Value = "'0.5947"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("D38").Value = "'2.582"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.88%  "
$ws.Range("D39").Value = "'0.8547"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").Value = "'0.9994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -7.58%  "
$ws.Range("D42").Value = "'98.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").Value = "'1.792"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.33%  "
$ws.Range("D44").Value = "'0.3686"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.90%  "
$ws.Range("D45").Value = "'4.643"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.54%  "
$ws.Range("D46").Value = "'0.1094"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.26%  "
$ws.Range("D49").Value = "'29.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.27%  "
$ws.Range("D50").Value = "'0.9994"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  +0.27%  "

# Rows 47 and 48 swap content: Cronos <-> Aptos (row 47 becomes Aptos, row 48 becomes Cronos)
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'6.101"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.64%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05220"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "
